$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold font, borders, alignment) from the existing "sum"
# header cell so the new "Save" header matches the other header cells'
# style (same shared cellXf), then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data column: 1 = saved this row, 0 = not saved.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
